$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists forms alphabetically by name (column A) with a matching
# URL in column B. "Civil No Contact Order - CNCO" sorts between
# "Cannabis expungement" (row 4) and "Criminal Court fee waiver" (row 5),
# so insert a new row 5 for it and push everything else down by one.

# Remember, in original row order, which rows currently carry a hyperlink
# so they can be rebuilt (shifted) after the insert - this engine's
# Rows.Insert() does not renumber the Hyperlinks collection automatically.
$oldHyperlinkRows = @(2,16,13,40,18,44,34,9,10,11,36,47,19,21,45,30,31,29,28,25,23,4,12,41,33,38,42,39,5,17,35,46,22,6)

# Drop all existing hyperlink objects (cell text/format is untouched); they
# get recreated below at their shifted locations.
$ws.Hyperlinks.Delete()

$ws.Rows("5:5").Insert()

$ws.Cells.Item(5, 1).Value = "Civil No Contact Order - CNCO"
$ws.Cells.Item(5, 2).Value = "https://www.illinoislegalaid.org/legal-information/civil-no-contact-order-request"
$ws.Cells.Item(5, 2).Style = "Hyperlink"

foreach ($oldRow in $oldHyperlinkRows) {
    $newRow = $oldRow
    if ($oldRow -ge 5) {
        $newRow = $oldRow + 1
    }
    $cell = $ws.Cells.Item($newRow, 2)
    $ws.Hyperlinks.Add($cell, $cell.Text) | Out-Null
    # Hyperlinks.Add() re-applies hyperlink formatting through a fresh
    # style slot; reassert the normal "Hyperlink" cell style so the cell
    # keeps using the workbook's existing (shared) style record.
    $cell.Style = "Hyperlink"
}

# The worksheet keeps a record of the last sort (rows 2-16 originally);
# extend it to rows 2-17 now that the new row pushed everything down.
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("A2:A17"))
$so.SetRange($ws.Range("A2:B17"))
$so.Header = 0
$so.Apply()

$ws.Range("B8").Select() | Out-Null
